$wb = $excel.ActiveWorkbook

# --- Insert a new "item_num" column (column O) on every sheet that shares
#     the common Measures-style header layout (A:O with "comment" in O).
#     The existing column O (comment) shifts right to become column P.
$sheetNames = @("Measures", "ID", "Dems", "Dates", "NewVars")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns("O").Insert()
    $ws.Range("O1").Value = "item_num"
}

# --- Fill in item_num values on the Measures sheet (rows 2-49 = 1, rows 50-51 = 4)
$measures = $wb.Worksheets.Item("Measures")
for ($r = 2; $r -le 49; $r++) {
    $measures.Cells.Item($r, 15).Value = 1
}
for ($r = 50; $r -le 51; $r++) {
    $measures.Cells.Item($r, 15).Value = 4
}

# --- Update the _FilterDatabase defined name so it covers the new column
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Measures!_FilterDatabase") {
        $n.RefersTo = "=Measures!`$A`$1:`$P`$51"
    }
}

# --- Fix up selections on the sheets whose "I1:I1048576" selection shifted to "O1:O1048576"
foreach ($name in @("ID", "Dems", "Dates", "NewVars")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("O1:O1048576").Select()
}

# --- Make Measures the active sheet/tab, with O2 selected (inside the frozen pane)
$measures.Activate()
$measures.Range("O2").Select()
